$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hats")
$ws.Activate()

# Add "x" markers in column A for the hats that were added (issue #58)
$rows = @(19, 20, 21, 22, 23, 24, 36, 43, 44, 48)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = "x"
}

# Leave the view scrolled/selected where editing left off
$ws.Range("A28").Select()
$ws.Range("A49").Select()
$excel.ActiveWindow.ScrollRow = 28
